# Reverts the presentation to a previous version:
#  - slide 1: bumps the "Antigravity Agent" date
#  - slide 5: updates several model-comparison table metrics
#  - slide 6: updates the LightGBM convergence hyper-parameters

$p = $ppt.ActivePresentation

# --- Slide 1: subtitle date line -------------------------------------------
$s1 = $p.Slides.Item(1)
$dateRun = $s1.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).TextRange.Runs(1)
$dateRun.Text = "Antigravity Agent – 2025-12-11"

# --- Slide 5: "Comparaison Globale des Résultats" table ---------------------
$s5 = $p.Slides.Item(5)
$tbl = $s5.Shapes.Item(3).Table

# Row 2 = LightGBM (Vainqueur)
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Runs(1).Text = "0.790"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Runs(1).Text = "0.835"
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Runs(1).Text = "72.3%"

# Row 3 = GRU (RNN)
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Runs(1).Text = "0.713"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Runs(1).Text = "0.743"
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Runs(1).Text = "65.1%"
$tbl.Cell(3, 5).Shape.TextFrame.TextRange.Runs(1).Text = "0.217"

# Row 4 = Transformer
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Runs(1).Text = "0.711"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Runs(1).Text = "0.748"
$tbl.Cell(4, 4).Shape.TextFrame.TextRange.Runs(1).Text = "66.3%"
$tbl.Cell(4, 5).Shape.TextFrame.TextRange.Runs(1).Text = "0.211"

# --- Slide 6: convergence bullet on the LightGBM (Optuna) text box ---------
$s6 = $p.Slides.Item(6)
$convBox = $s6.Shapes.Item(3)
$origHeight = $convBox.Height
$convRun = $convBox.TextFrame.TextRange.Paragraphs(4).TextRange.Runs(1)
$convRun.Text = "•  Convergence vers paramètres robustes (n_est=318, lr=0.018)."
# The textbox has spAutoFit; keep its original height (only the wording changed).
$convBox.Height = $origHeight
